# Applies the "Login" sheet addition + new validation strings to the
# ValidationStrings workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Create Account" sheet (sheet2): new rows 2-4 (Resend code / invalid
#    code / invalid email scenarios).
# ---------------------------------------------------------------------
$wsCreateAccount = $wb.Worksheets.Item("Create Account")

$wsCreateAccount.Range("A2").Value = "ResendCodeMessage"
$wsCreateAccount.Range("B2").Value = "//div[@class='alert__title']"
$wsCreateAccount.Range("C2").Value = "確認コードを再送しました、宛先メールをご確認ください。"

$wsCreateAccount.Range("A3").Value = "InvalidCodeError"
$wsCreateAccount.Range("B3").Value = "//p[@class='alert__des']"
$wsCreateAccount.Range("C3").Value = "6桁の確認コードを入力してください。"

$wsCreateAccount.Range("A4").Value = "InvalidEmailErrorMessage"
$wsCreateAccount.Range("B4").Value = "//p[@class='alert__des']"
$wsCreateAccount.Range("C4").Value = "メールアドレスは、メールアドレス形式で入力してください。"

# ---------------------------------------------------------------------
# 2) "Reset password" sheet (sheet1): new row 5 (forgot-password message),
#    with a small purple Arial font on the expected-string cell.
# ---------------------------------------------------------------------
$wsResetPassword = $wb.Worksheets.Item("Reset password")

$wsResetPassword.Range("A5").Value = "ForgotPasswordPageMessage2"
$wsResetPassword.Range("B5").Value = "//p[@class='ep-reset-pass__desc']"
$wsResetPassword.Range("C5").Value = "パスワード再設定のメールをお送りします。"

$c5Font = $wsResetPassword.Range("C5").Font
$c5Font.Size = 7
$c5Font.Color = 6563904
$c5Font.Name = "Arial"

# ---------------------------------------------------------------------
# 3) Add the new "Login" sheet at the end, with its own validation rows.
# ---------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$wsLogin = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($sheetCount))
$wsLogin.Name = "Login"

$wsLogin.Range("A1").Value = "Label"
$wsLogin.Range("B1").Value = "Locator"
$wsLogin.Range("C1").Value = "ExpectedString"
$wsLogin.Range("A1:C1").Font.Bold = $true

$wsLogin.Range("A2").Value = "InvalidDetailsError"
$wsLogin.Range("B2").Value = "//p[@class='alert__des']"
$wsLogin.Range("C2").Value = "メールアドレスは、メールアドレス形式で入力してください。"

$wsLogin.Range("A3").Value = "InvalidUserNameError"
$wsLogin.Range("B3").Value = "//p[@class='alert__des']"
$wsLogin.Range("C3").Value = "メールアドレスは、メールアドレス形式で入力してください。"

$wsLogin.Range("A4").Value = "InvalidPasswordError"
$wsLogin.Range("B4").Value = "//p[@class='alert__des']"
$wsLogin.Range("C4").Value = "メールアドレスまたEiDまたはパスワードをご確認してください。"

# ---------------------------------------------------------------------
# 4) Selections: restore/update the active cell on each sheet, finishing
#    on "Login" so it ends up the active (tabSelected) tab.
# ---------------------------------------------------------------------
$wsCreateAccount.Range("B4").Select()
$wsResetPassword.Range("C5").Select()
$wsLogin.Range("B3:B4").Select()
